# Update the "想去人数" (want-to-go count) figures that changed between
# the two gh-pages data generations.
#
# Sheet "展览" (Exhibitions):
#   Row 4 (南宁·万圣漫控嘉年华10)  F4: 833 -> 834
#   Row 6 (南宁·黑塔利亚同人ONLY)  F6: 26  -> 27
#
# Sheet "全部类型" (All types):
#   Row 5 (南宁·万圣漫控嘉年华10)  F5: 833 -> 834
#   Row 7 (南宁·黑塔利亚同人ONLY)  F7: 26  -> 27

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 834
$wsExhibition.Range("F6").Value = 27

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F5").Value = 834
$wsAllTypes.Range("F7").Value = 27
